$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-34 (A: id, B: business_status, C: name, D: rating, E: user_ratings_total)
$data = @(
    @(8,'OPERATIONAL','500 Festival',4.5,10),
    @(31,'OPERATIONAL','Arts for Lawrence',4.9,19),
    @(7,'OPERATIONAL','Carmel International Arts Festival',4.8,18),
    @(18,'OPERATIONAL','Celebration Plaza Amphitheater',4.6,26),
    @(20,'OPERATIONAL','Cumberland Arts, Inc.',0,0),
    @(21,'OPERATIONAL','Downtown Indy, Inc.',4.6,246),
    @(19,'OPERATIONAL','Festival Flea Market',3.9,191),
    @(23,'OPERATIONAL','Free Stage',4.6,361),
    @(22,'OPERATIONAL','Greenfield Main Street, Inc.',5,4),
    @(13,'OPERATIONAL','HSI Show Productions',0,0),
    @(17,'OPERATIONAL','Heartland Film',4.3,14),
    @(25,'OPERATIONAL','Holy Rosary Church',4.8,210),
    @(16,'OPERATIONAL','Indiana Arts Building',4,3),
    @(12,'OPERATIONAL','Indianapolis Auto Show',3.1,51),
    @(14,'OPERATIONAL','Indianapolis Jewish Film Festival',0,0),
    @(27,'OPERATIONAL','Indianapolis Zoo',4.5,14982),
    @(0,'OPERATIONAL','Indy Taco Fest',4.2,30),
    @(30,'OPERATIONAL','IndyFringe Theatre',4.6,170),
    @(1,'OPERATIONAL','Irvington Halloween Festival',4.7,13),
    @(32,'OPERATIONAL','Lights Over Seminary',5,2),
    @(11,'OPERATIONAL','Lourdes Lyons Chili Cook-Off',0,0),
    @(15,'OPERATIONAL','MasterWorks Festival',0,0),
    @(10,'OPERATIONAL','Military Park',4.6,1512),
    @(29,'OPERATIONAL','Nickel Plate District Amphitheater',4.6,430),
    @(24,'OPERATIONAL','Pendleton Fall Creek Heritage Fair',0,0),
    @(6,'OPERATIONAL','Riley Festival',4.3,71),
    @(9,'OPERATIONAL','Spirit & Place Festival',0,0),
    @(3,'OPERATIONAL','St Louis de Montfort Fall Festival',4,4),
    @(4,'OPERATIONAL','St. Christopher Mid-Summer Festival',4.5,13),
    @(26,'OPERATIONAL','Stonycreek Farm Presents: HauntFest',4,71),
    @(5,'OPERATIONAL','Talbot Street Art Fair',4.5,99),
    @(28,'OPERATIONAL','Wild Blackberry Farms',5,81),
    @(2,'OPERATIONAL','indianapolis greek fest',4.3,19)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

# Remove now-extra rows (old sheet had data through row 38, new one ends at row 34)
$ws.Range("A35:E38").Delete()
